# Fix: resolve timetable regeneration issue with file uploads
# Updates course codes/instructors (old 1xx codes -> new 5x codes) across
# Section_A, Section_B timetable grids and the Course_Summary sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Section_A
# ---------------------------------------------------------------
$wsA = $wb.Worksheets.Item("Section_A")

$wsA.Range("B2").Value = "EC105"
$wsA.Range("C2").Value = "Free"
$wsA.Range("D2").Value = "Free"
$wsA.Range("E2").Value = "CS105"
$wsA.Range("F2").Value = "MA105"

$wsA.Range("B3").Value = "Free"
$wsA.Range("C3").Value = "Free"
$wsA.Range("D3").Value = "HS105"
$wsA.Range("E3").Value = "DS105"
$wsA.Range("F3").Value = "CS105"

$wsA.Range("B5").Value = "MA105"
$wsA.Range("C5").Value = "EC105"
$wsA.Range("D5").Value = "Free"
$wsA.Range("E5").Value = "Free"
$wsA.Range("F5").Value = "MA106"

$wsA.Range("B6").Value = "CS105"
$wsA.Range("C6").Value = "HS105"
$wsA.Range("D6").Value = "Free"
$wsA.Range("E6").Value = "CS155 (Elective)"
$wsA.Range("F6").Value = "EC105"

$wsA.Range("B7").Value = "HS105"
$wsA.Range("C7").Value = "DS105"
$wsA.Range("D7").Value = "MA106"
$wsA.Range("E7").Value = "Free"

# ---------------------------------------------------------------
# Section_B
# ---------------------------------------------------------------
$wsB = $wb.Worksheets.Item("Section_B")

$wsB.Range("B2").Value = "CS105"
$wsB.Range("C2").Value = "Free"
$wsB.Range("D2").Value = "Free"
$wsB.Range("E2").Value = "DS105"
$wsB.Range("F2").Value = "MA105"

$wsB.Range("C3").Value = "MA105"
$wsB.Range("D3").Value = "CS105"
$wsB.Range("E3").Value = "EC105"
$wsB.Range("F3").Value = "Free"

$wsB.Range("C5").Value = "HS105"
$wsB.Range("D5").Value = "HS105"
$wsB.Range("E5").Value = "HS105"
$wsB.Range("F5").Value = "CS105"

$wsB.Range("B6").Value = "Free"
$wsB.Range("C6").Value = "MA106"
$wsB.Range("D6").Value = "DS105"
$wsB.Range("E6").Value = "CS155 (Elective)"

$wsB.Range("B7").Value = "EC105"
$wsB.Range("C7").Value = "Free"
$wsB.Range("E7").Value = "MA106"
$wsB.Range("F7").Value = "EC105"

# ---------------------------------------------------------------
# Course_Summary
# ---------------------------------------------------------------
$wsC = $wb.Worksheets.Item("Course_Summary")

$wsC.Range("A2").Value = "MA105"
$wsC.Range("B2").Value = "Statistical Methods"
$wsC.Range("F2").Value = "Dr. Priya Sharma"

$wsC.Range("A3").Value = "DS105"
$wsC.Range("B3").Value = "Fundamentals of Data Science"
$wsC.Range("F3").Value = "Dr. Rohan Verma"

$wsC.Range("A4").Value = "MA106"
$wsC.Range("B4").Value = "Probability Theory"
$wsC.Range("F4").Value = "Dr. Sanjay Kumar"

$wsC.Range("A5").Value = "EC105"
$wsC.Range("B5").Value = "Digital Systems"
$wsC.Range("F5").Value = "Dr. Anjali Mehta"

$wsC.Range("A6").Value = "CS105"
$wsC.Range("B6").Value = "Programming Fundamentals"
$wsC.Range("F6").Value = "Dr. Vikram Joshi"

$wsC.Range("A7").Value = "HS105"
$wsC.Range("B7").Value = "Professional Communication"
$wsC.Range("F7").Value = "Dr. Rajeev Malhotra"

$wsC.Range("A8").Value = "CS155"
$wsC.Range("B8").Value = "Cybersecurity Basics"
$wsC.Range("F8").Value = "Dr. Kavya Iyer"

$wb.Save()
